$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.055.22"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.312.42"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "2.673.20"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "2.330.56"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").Value = "42.993.42"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.10%  "
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E29").Value = "  -10.39%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.05%  "
$ws.Range("E33").Value = "  +3.51%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.109"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").Value = "1.993.96"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.94%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").Value = "2.539.27"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  -0.21%  "
